$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new cells are treated as plain text so numeric-looking
# strings (IDs, "1.599.00", "83.90") are stored verbatim, matching the
# inlineStr cells already used throughout the sheet.
$ws.Range("A107:A108").NumberFormat = "@"
$ws.Range("C107:C108").NumberFormat = "@"

$ws.Range("A107").Value = "89975816"
$ws.Range("B107").Value = "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Clean Deca"
$ws.Range("C107").Value = "1.599.00"

$ws.Range("A108").Value = "90308603"
$ws.Range("B108").Value = "Torneira Elétrica Bica Alta Branca 220V 5500W Prime Equation"
$ws.Range("C108").Value = "83.90"
